$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three rows (5:7) - the new TPM data only has 3 data rows
$ws.Rows("5:7").Delete()

# Row 2 - new values (MuSCs / Bdnf / Ntrk2 / ECs)
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Bdnf"
$ws.Range("C2").Value = "Ntrk2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.9253576666666667
$ws.Range("H2").Value = 2.776073
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1433513333333333
$ws.Range("N2").Value = 0.430054
$ws.Range("O2").Value = 0.002710896760593916
$ws.Range("P2").Value = 0.002710896760593916
$ws.Range("Q2").Value = 0.1326512553268889
$ws.Range("R2").Value = 1.193861297942
$ws.Range("S2").Value = 0.002710896760593916
$ws.Range("T2").Value = 0.002710896760593916

# Row 3 - new values (MuSCs / Bdnf / Ntrk2 / FAPs)
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Bdnf"
$ws.Range("C3").Value = "Ntrk2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.9253576666666667
$ws.Range("H3").Value = 2.776073
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 43.12631833333334
$ws.Range("N3").Value = 129.378955
$ws.Range("O3").Value = 0.8155556976531461
$ws.Range("P3").Value = 0.8155556976531461
$ws.Range("Q3").Value = 39.90726930485723
$ws.Range("R3").Value = 359.1654237437151
$ws.Range("S3").Value = 0.8155556976531461
$ws.Range("T3").Value = 0.8155556976531461

# Row 4 - new values (MuSCs / Bdnf / Ntrk2 / MuSCs)
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Bdnf"
$ws.Range("C4").Value = "Ntrk2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.9253576666666667
$ws.Range("H4").Value = 2.776073
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.610002999999999
$ws.Range("N4").Value = 28.830009
$ws.Range("O4").Value = 0.18173340558626
$ws.Range("P4").Value = 0.1817334055862599
$ws.Range("Q4").Value = 8.892689952739666
$ws.Range("R4").Value = 80.034209574657
$ws.Range("S4").Value = 0.18173340558626
$ws.Range("T4").Value = 0.1817334055862599
